# Apply the FHIR R4B -> R4 regeneration edits to the StructureDefinition
# spreadsheet (study-family-focus extension).

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet -------------------------------------------------
# Row 8:  Date            A8="Date"          B8 = generation timestamp
# Row 15: FHIR Version    A15="FHIR Version" B15 = "4.3.0" -> "4.0.1"
$wsMeta.Range("B8").Value = "2025-06-13T15:45:04+00:00"
$wsMeta.Range("B15").Value = "4.0.1"

# --- Elements sheet ---------------------------------------------------
# Row 2 = Extension (root).  Constraint(s) column = AJ.
# The ele-1 invariant text loses the "unless an empty Parameters resource
# ... or $this is Parameters" clause, matching the plain R4 wording that
# already appears elsewhere (row 4, Element.extension).
$wsElem.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}
ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Row 3 = Extension.id.  Type(s) column = K: "id" -> "string"
$wsElem.Range("K3").Value = "string
"

# Row 6 = Extension.value[x]. Definition column = M: R4B doc link -> R4
$wsElem.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
